# Adds more League of Legends champion-mastery game data to the sheet:
# two new header columns (Hover / Party Bonus) and fifteen new data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns -----------------------------------------------
$ws.Range("K1").Value = "Hover"
$ws.Range("L1").Value = "Party Bonus"

# --- Row 10 -------------------------------------------------------------
$ws.Range("A10").Value = "Jazzrich"
$ws.Range("B10").Value = "Zac"
$ws.Range("F10").Value = "B+"
$ws.Range("G10").Value = $false
$ws.Range("J10").Value = 189

# --- Row 11 -------------------------------------------------------------
$ws.Range("A11").Value = "Keno"
$ws.Range("B11").Value = "Jhin"
$ws.Range("D11").Value = 1.5631944444444443
$ws.Range("D11").NumberFormat = "[h]:mm:ss"
$ws.Range("G11").Value = $false
$ws.Range("J11").Value = 229

# --- Row 12 -------------------------------------------------------------
$ws.Range("A12").Value = "Albert"
$ws.Range("B12").Value = "Blitz"
$ws.Range("D12").Value = 0.90069444444444446
$ws.Range("D12").NumberFormat = "h:mm"
$ws.Range("F12").Value = "B"
$ws.Range("G12").Value = $false
$ws.Range("J12").Value = 167

# --- Row 13 -------------------------------------------------------------
$ws.Range("A13").Value = "Jazzrich"
$ws.Range("B13").Value = "Maokai"
$ws.Range("D13").Value = 0.90069444444444446
$ws.Range("D13").NumberFormat = "h:mm"
$ws.Range("F13").Value = "C+"
$ws.Range("G13").Value = $false
$ws.Range("J13").Value = 167

# --- Row 14 -------------------------------------------------------------
$ws.Range("A14").Value = "NeoStrykr007"
$ws.Range("B14").Value = "Azir"
$ws.Range("D14").Value = 0.90069444444444446
$ws.Range("D14").NumberFormat = "h:mm"
$ws.Range("F14").Value = "B+"
$ws.Range("G14").Value = $false
$ws.Range("J14").Value = 167

# --- Row 15 -------------------------------------------------------------
$ws.Range("A15").Value = "Greatest Forever"
$ws.Range("B15").Value = "Lucian"
$ws.Range("D15").Value = 2.25
$ws.Range("D15").NumberFormat = "[h]:mm:ss"
$ws.Range("F15").Value = "S-"
$ws.Range("G15").Value = $false
$ws.Range("J15").Value = 337
$ws.Range("L15").Value = 11

# --- Row 16 -------------------------------------------------------------
$ws.Range("A16").Value = "Random"
$ws.Range("B16").Value = "Ezreal"
$ws.Range("D16").Value = 1.2395833333333333
$ws.Range("D16").NumberFormat = "[h]:mm:ss"
$ws.Range("F16").Value = "S-"
$ws.Range("G16").Value = $true
$ws.Range("K16").Value = 220
$ws.Range("L16").Value = 11

# --- Row 17 -------------------------------------------------------------
$ws.Range("A17").Value = "Seo"
$ws.Range("B17").Value = "Riven"
$ws.Range("D17").Value = 1.2395833333333333
$ws.Range("D17").NumberFormat = "[h]:mm:ss"
$ws.Range("F17").Value = "B"
$ws.Range("G17").Value = $false
$ws.Range("J17").Value = 197
$ws.Range("K17").Value = 35
$ws.Range("L17").Value = 11

# --- Row 18 -------------------------------------------------------------
$ws.Range("A18").Value = "Greatest Forever"
$ws.Range("B18").Value = "Kalista"
$ws.Range("D18").Value = 1.2395833333333333
$ws.Range("D18").NumberFormat = "[h]:mm:ss"
$ws.Range("F18").Value = "B+"
$ws.Range("G18").Value = $false
$ws.Range("J18").Value = 197
$ws.Range("K18").Value = 37
$ws.Range("L18").Value = 11

# --- Row 19 -------------------------------------------------------------
$ws.Range("A19").Value = "Ryue"
$ws.Range("B19").Value = "Kha'Zix"
$ws.Range("D19").Value = 1.2395833333333333
$ws.Range("D19").NumberFormat = "[h]:mm:ss"
$ws.Range("F19").Value = "A"
$ws.Range("G19").Value = $false
$ws.Range("J19").Value = 197
$ws.Range("L19").Value = 11

# --- Row 20 -------------------------------------------------------------
$ws.Range("A20").Value = "Jazzrich"
$ws.Range("B20").Value = "Draven"
$ws.Range("D20").Value = 1.7388888888888889
$ws.Range("D20").NumberFormat = "[h]:mm:ss"
$ws.Range("F20").Value = "A+"
$ws.Range("G20").Value = $true
$ws.Range("J20").Value = 1458
$ws.Range("K20").Value = 278
$ws.Range("L20").Value = 82

# --- Row 21 -------------------------------------------------------------
$ws.Range("A21").Value = "NeoStrykr007"
$ws.Range("B21").Value = "Sion"
$ws.Range("D21").Value = 1.7388888888888889
$ws.Range("D21").NumberFormat = "[h]:mm:ss"
$ws.Range("F21").Value = "B+"
$ws.Range("G21").Value = $true
$ws.Range("J21").Value = 1458
$ws.Range("K21").Value = 254
$ws.Range("L21").Value = 82

# --- Row 22 -------------------------------------------------------------
$ws.Range("A22").Value = "Albert"
$ws.Range("B22").Value = "Leona"
$ws.Range("D22").Value = 1.7388888888888889
$ws.Range("D22").NumberFormat = "[h]:mm:ss"
$ws.Range("F22").Value = "A+"
$ws.Range("G22").Value = $true
$ws.Range("J22").Value = 1458
$ws.Range("K22").Value = 283
$ws.Range("L22").Value = 82

# --- Row 23 -------------------------------------------------------------
$ws.Range("A23").Value = "?"
$ws.Range("B23").Value = "?"
$ws.Range("D23").Value = 1.7659722222222223
$ws.Range("D23").NumberFormat = "[h]:mm:ss"
$ws.Range("F23").Value = "A+"
$ws.Range("G23").Value = $true
$ws.Range("J23").Value = 1428
$ws.Range("K23").Value = 281
$ws.Range("L23").Value = 54

# --- Row 24 -------------------------------------------------------------
$ws.Range("A24").Value = "?"
$ws.Range("B24").Value = "?"
$ws.Range("D24").Value = 1.7659722222222223
$ws.Range("D24").NumberFormat = "[h]:mm:ss"
$ws.Range("F24").Value = "B+"
$ws.Range("G24").Value = $true
$ws.Range("J24").Value = 1428
$ws.Range("K24").Value = 287
$ws.Range("L24").Value = 54

# --- Column A got wider to fit the longest new player name --------------
$ws.Columns.Item(1).ColumnWidth = 15.142857142857142

# --- Selection moved as part of editing the sheet ------------------------
$ws.Range("E17").Select() | Out-Null
